# edit.ps1 - Apply "Updated cryptos list" data refresh (Wed Jul 5 14:47:36 UTC 2023)
#
# Updates Price (D) and Volume(1h) (E) columns for most rows with refreshed
# market data, and additionally re-sorts two adjacent pairs of coins whose rank
# order flipped (BitcoinCash/WrappedBTC at rows 16-17, and Elrond/Algorand at
# rows 49-50) by swapping their Coin/Link/Price/Volume cells in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores values that look numeric but must stay literal TEXT
# (e.g. multi-dot thousands separators like "30.421.40", or trailing zeros that
# a real number would drop). Force text formatting before writing so Excel does
# not coerce these into numeric cells.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.421.40'
$ws.Range("E2").Value = '  -1.97%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.909.63'
$ws.Range("E3").Value = '  -2.55%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  -0.20%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.54'
$ws.Range("E5").Value = '  -2.50%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("E6").Value = '  -0.16%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4760'
$ws.Range("E7").Value = '  -2.30%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2862'
$ws.Range("E8").Value = '  -3.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06702'
$ws.Range("E9").Value = '  -3.65%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.82'
$ws.Range("E10").Value = '  -3.56%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '102.85'
$ws.Range("E11").Value = '  -4.69%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07716'
$ws.Range("E12").Value = '  -1.14%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.907.62'
$ws.Range("E13").Value = '  -3.97%  '

# Row 14
$ws.Range("E14").Value = '  -5.31%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6772'
$ws.Range("E15").Value = '  -3.35%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9997'
$ws.Range("E18").Value = '  -0.06%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007507'
$ws.Range("E19").Value = '  -3.38%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.74'
$ws.Range("E20").Value = '  -4.36%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.419'
$ws.Range("E21").Value = '  -2.18%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("E22").Value = '  -0.18%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.314'
$ws.Range("E23").Value = '  -3.13%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.468'
$ws.Range("E24").Value = '  -4.01%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '164.52'
$ws.Range("E25").Value = '  -2.47%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.00'
$ws.Range("E26").Value = '  -4.91%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.060'
$ws.Range("E27").Value = '  -6.04%  '

# Row 28
$ws.Range("E28").Value = '  -3.73%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.377'
$ws.Range("E29").Value = '  -0.85%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.619'
$ws.Range("E30").Value = '  -0.47%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.514'
$ws.Range("E31").Value = '  -3.60%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.262'
$ws.Range("E32").Value = '  -4.71%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04790'
$ws.Range("E33").Value = '  -2.65%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7320'
$ws.Range("E34").Value = '  -3.18%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.116'
$ws.Range("E35").Value = '  -4.85%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9986'
$ws.Range("E36").Value = '  -0.21%  '

# Row 38
$ws.Range("E38").Value = '  -4.10%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.601'
$ws.Range("E39").Value = '  -3.91%  '

# Row 40
$ws.Range("E40").Value = '  -4.86%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '75.00'
$ws.Range("E41").Value = '  -3.89%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.993'
$ws.Range("E42").Value = '  -6.65%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8674'
$ws.Range("E43").Value = '  -3.93%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '106.73'
$ws.Range("E44").Value = '  -2.64%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.058.05'
$ws.Range("E45").Value = '  +4.41%  '

# Row 46
$ws.Range("E46").Value = '  -4.29%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9989'
$ws.Range("E47").Value = '  -0.22%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.493'
$ws.Range("E48").Value = '  -7.71%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.897'
$ws.Range("E51").Value = '  -5.66%  '

# Rows 16-17: BitcoinCash and WrappedBTC swap places in the ranking
# Row 16
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.439.56'
$ws.Range("E16").Value = '  -1.95%  '

# Row 17
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '259.80'
$ws.Range("E17").Value = '  -7.62%  '

# Rows 49-50: Elrond and Algorand swap places in the ranking
# Row 49
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1202'
$ws.Range("E49").Value = '  -4.30%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.15'
$ws.Range("E50").Value = '  -2.37%  '

